# "actualizando la tabla de tareas y escogiendo algunas"
#
# Fill in the "Buscar una base de datos con alimentos" (row 16) and
# "Volcar una bd con ingredientes en la nuestra / Agregarlos a mano"
# (row 17) sub-task rows of the Tareas table: assign the responsible
# person (Mario), record estimated/invested time, and mark row 16 as
# done with the same green check fill used elsewhere in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: "Buscar una base de datos con alimentos" -> done, 1 hora / 1 hora
$ws.Range("C16").Interior.Color = 5287936
$ws.Range("D16").Value = "Mario"
$ws.Range("E16").Value = "1 hora"
$ws.Range("F16").Value = "1 hora"

# Row 17: "Volcar una bd con ingredientes en la nuestra / Agregarlos a mano"
$ws.Range("D17").Value = "Mario"
$ws.Range("E17").Value = "20 min"

# Match the author's final on-screen selection/position.
$ws.Range("C16").Select()
